$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.438.36'
$ws.Range("E2").Value = '  +0.10%  '

$ws.Range("D3").Value = '2.719.75'
$ws.Range("E3").Value = '  +2.82%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '''608.79'
$ws.Range("E5").Value = '  +1.99%  '

$ws.Range("D6").Value = '''167.54'
$ws.Range("E6").Value = '  +5.66%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("E8").Value = '  +2.45%  '

$ws.Range("D9").Value = '2.717.58'
$ws.Range("E9").Value = '  +2.74%  '

$ws.Range("D10").Value = '''0.145'
$ws.Range("E10").Value = '  +1.83%  '

$ws.Range("E11").Value = '  +4.49%  '

$ws.Range("E13").Value = '  +0.62%  '

$ws.Range("D14").Value = '''28.58'
$ws.Range("E14").Value = '  +2.27%  '

$ws.Range("D15").Value = '3.217.82'
$ws.Range("E15").Value = '  +3.00%  '

$ws.Range("E16").Value = '  +1.19%  '

$ws.Range("D17").Value = '68.394.54'
$ws.Range("E17").Value = '  +0.29%  '

$ws.Range("D18").Value = '2.723.93'
$ws.Range("E18").Value = '  +2.61%  '

$ws.Range("D19").Value = '''11.89'
$ws.Range("E19").Value = '  +4.31%  '

$ws.Range("D20").Value = '''372.26'
$ws.Range("E20").Value = '  +2.51%  '

$ws.Range("D21").Value = '''7.65'
$ws.Range("E21").Value = '  +2.62%  '

$ws.Range("D22").Value = '''4.51'
$ws.Range("E22").Value = '  +2.68%  '

$ws.Range("D23").Value = '''4.99'
$ws.Range("E23").Value = '  +4.81%  '

$ws.Range("D24").Value = '''2.09'
$ws.Range("E24").Value = '  +1.41%  '

$ws.Range("D25").Value = '''73.12'
$ws.Range("E25").Value = '  -1.74%  '

$ws.Range("E26").Value = '  +0.00%  '

$ws.Range("D27").Value = '''10.12'
$ws.Range("E27").Value = '  +3.86%  '

$ws.Range("D29").Value = '''0.0000104'
$ws.Range("E29").Value = '  +1.42%  '

$ws.Range("D30").Value = '''587.28'
$ws.Range("E30").Value = '  +5.03%  '

$ws.Range("D31").Value = '''1.00'
$ws.Range("E31").Value = '  -0.13%  '

$ws.Range("D32").Value = '''8.20'
$ws.Range("E32").Value = '  +1.95%  '

$ws.Range("D33").Value = '''1.43'
$ws.Range("E33").Value = '  +3.12%  '

$ws.Range("E34").Value = '  +6.65%  '

$ws.Range("E35").Value = '  +2.09%  '

$ws.Range("D36").Value = '''1.60'
$ws.Range("E36").Value = '  -2.86%  '

$ws.Range("D38").Value = '''162.72'
$ws.Range("E38").Value = '  +2.01%  '

$ws.Range("D39").Value = '''19.90'
$ws.Range("E39").Value = '  +1.41%  '

$ws.Range("D40").Value = '''0.378'
$ws.Range("E40").Value = '  +2.36%  '

$ws.Range("D41").Value = '''1.89'
$ws.Range("E41").Value = '  +1.33%  '

$ws.Range("D42").Value = '''5.43'
$ws.Range("E42").Value = '  +2.34%  '

$ws.Range("D43").Value = '''17.98'
$ws.Range("E43").Value = '  +0.97%  '

$ws.Range("E44").Value = '  +1.19%  '

$ws.Range("E45").Value = '  +0.02%  '

$ws.Range("D46").Value = '0.0₆0310'
$ws.Range("E46").Value = '  -3.24%  '

$ws.Range("D47").Value = '''40.98'
$ws.Range("E47").Value = '  +1.72%  '

$ws.Range("D48").Value = '''0.598'
$ws.Range("E48").Value = '  +4.32%  '

$ws.Range("D49").Value = '''155.00'
$ws.Range("E49").Value = '  -1.97%  '

$ws.Range("E50").Value = '  +3.51%  '

$ws.Range("E51").Value = '  +5.23%  '
